$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 194, shifting existing rows 194-290 down to 195-291.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row 194 with the new record.
$ws.Range("A194").Value = 9
$ws.Range("B194").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C194").Value = "Metropolitana"
$ws.Range("D194").Value = 44572
$ws.Range("D194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E194").Value = 13
$ws.Range("F194").Value = 100112032
$ws.Range("G194").Value = "Zapallo italiano"
$ws.Range("H194").Value = "Sin especificar"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 106
$ws.Range("K194").Value = 13000
$ws.Range("L194").Value = 15000
$ws.Range("M194").Value = 14000
$ws.Range("N194").Value = "`$/caja 50 unidades"
$ws.Range("O194").Value = "Región de O'Higgins"
$ws.Range("P194").Value = 280
$ws.Range("Q194").Value = 50
$ws.Range("R194").Value = "Hortaliza"
